# Generate Report for Handoff
#
# - Flips the "Status" column from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it appears (Overview!E2/F2, zh-cn!C2,
#   de-de!C2).
# - Refreshes the "Latest Handoff Datetime" / "Latest HO Xliff Generate
#   Date" timestamps to the new handoff run time.
# - Narrows the (now shorter) status column on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$readyStatus = "Ready for handoff"

# --- Status cells -----------------------------------------------------
$ws1.Range("E2").Value = $readyStatus
$ws1.Range("F2").Value = $readyStatus
$ws2.Range("C2").Value = $readyStatus
$ws3.Range("C2").Value = $readyStatus

# --- Timestamps ---------------------------------------------------------
# de-de handoff run (also mirrored on the Overview tab for this row)
$ws1.Range("G2").Value = "2016-08-26 10:57:04"
$ws3.Range("H2").Value = "2016-08-26 10:57:04"

# zh-cn handoff run
$ws2.Range("H2").Value = "2016-08-26 10:56:57"

# --- Column widths --------------------------------------------------
# Status column got shorter text, so narrow it on every sheet that shows
# it. The ColumnWidth setter quantizes internally, so feed it the inverse
# of the target stored width (17.2159881591797 chars) to land on the
# closest reachable value.
$newStatusColWidth = 16.38265482584637

$ws1.Range("E1").ColumnWidth = $newStatusColWidth
$ws1.Range("F1").ColumnWidth = $newStatusColWidth
$ws2.Range("C1").ColumnWidth = $newStatusColWidth
$ws3.Range("C1").ColumnWidth = $newStatusColWidth
